# Edit slide 9 ("Battle-Field-3"): the "Design patterns" bullet gains
# "singleton" to the pattern list and drops the spell-check "err" marks
# that used to sit on "fascade" and "iterator" (they become plain text
# merged into the leading run), matching:
#   "Design patterns – factory, memento, observer, fascade, iterator, "
#   "singleton, etc"
#   "."

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Find the target paragraph (the "Design patterns ..." bullet).
$paragraphs = $tr.Paragraphs()
for ($i = 1; $i -le $paragraphs.Count; $i++) {
    $candidate = $tr.Paragraphs($i)
    if ($candidate.Text -like "Design patterns*") {
        $para = $candidate
        break
    }
}

$prStart = $para.Start
$full = $para.Text

# 1) Replace the trailing "etc." with "singleton, etc." in one shot so the
#    new text takes on the formatting of that trailing (non-err) run.
$tailOffset = $full.IndexOf("etc.")
$tailRange = $tr.Characters($prStart + $tailOffset, 4)
$tailRange.Text = "singleton, etc."

# 2) Split the trailing "." into its own run by re-assigning it in place.
$full = $para.Text
$dotRange = $tr.Characters($prStart + $full.Length - 1, 1)
$dotRange.Text = "."

# 3) Merge everything up to and including "iterator, " into a single run
#    so the old "err=1" spell-check flags on "fascade" / "iterator" are
#    dropped (the merged run takes the formatting of its first run, which
#    has no err flag).
$mergeLen = $full.IndexOf("singleton") - 0
$mergeRange = $tr.Characters($prStart, $mergeLen)
$mergeRange.Text = $mergeRange.Text
